$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 through 499) holds the "Förändrad" date (serial 45172 -> 45175)
for ($row = 2; $row -le 499; $row++) {
    $ws.Cells.Item($row, 3).Value = 45175
}
